$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M2").Value = -293
$ws.Range("H2").Value = 865
$ws.Range("I2").Value = 406
$ws.Range("K2").Value = 406
$ws.Range("M12").Value = -205.75
$ws.Range("H12").Value = 714
$ws.Range("K12").Value = 375.75
$ws.Range("I12").Value = 375.75
$ws.Range("J12").Value = 1165
$ws.Range("L12").Value = 1165
$ws.Range("N12").Value = -1505
$ws.Range("L17").Value = 59602.962
$ws.Range("J17").Value = 19867.654
$ws.Range("H17").Value = 19867.654
$ws.Range("N17").Value = -59938.962
$ws.Range("I19").Value = 1375.25
$ws.Range("J19").Value = 25000972
$ws.Range("L19").Value = 25000972
$ws.Range("H19").Value = 12501174
$ws.Range("M19").Value = -1200.25
$ws.Range("K19").Value = 1375.25
$ws.Range("N19").Value = -25001322
$ws.Range("H33").Value = 581.3889
$ws.Range("L33").Value = 956.4286
$ws.Range("I33").Value = 342.72726
$ws.Range("N33").Value = -1414.4286
$ws.Range("J33").Value = 956.4286
$ws.Range("M33").Value = -113.72726
$ws.Range("K33").Value = 342.72726
$ws.Range("N86").Value = -9695.1113
$ws.Range("J86").Value = 7449.1113
$ws.Range("L86").Value = 7449.1113
$ws.Range("K86").Value = 1978.8667
$ws.Range("M86").Value = -855.8667
$ws.Range("H86").Value = 4030.2083
$ws.Range("I86").Value = 1978.8667
$ws.Range("K89").Value = 9894.333500000001
$ws.Range("I89").Value = 1978.8667
$ws.Range("L89").Value = 37245.5565
$ws.Range("H89").Value = 4030.2083
$ws.Range("N89").Value = -48477.5565
$ws.Range("M89").Value = -4278.333500000001
$ws.Range("J89").Value = 7449.1113
$ws.Range("J99").Value = 5991.1665
$ws.Range("K99").Value = 521.25
$ws.Range("I99").Value = 173.75
$ws.Range("M99").Value = 976.75
$ws.Range("H99").Value = 3664.2
$ws.Range("N99").Value = -20969.4995
$ws.Range("L99").Value = 17973.4995
$ws.Range("I132").Value = 4501.375
$ws.Range("H132").Value = 8569.950000000001
$ws.Range("M132").Value = -10974.125
$ws.Range("J132").Value = 11282.333
$ws.Range("K132").Value = 13504.125
$ws.Range("N132").Value = -38906.999
$ws.Range("L132").Value = 33846.999
$ws.Range("L135").Value = 64914.4296
$ws.Range("J135").Value = 7212.7144
$ws.Range("I135").Value = 469.96667
$ws.Range("K135").Value = 4229.70003
$ws.Range("H135").Value = 1745.6216
$ws.Range("N135").Value = -69984.4296
$ws.Range("M135").Value = -1694.70003
$ws.Range("J138").Value = 5098.909
$ws.Range("M138").Value = 1195.3216
$ws.Range("K138").Value = 3944.6784
$ws.Range("L138").Value = 15296.727
$ws.Range("N138").Value = -25576.727
$ws.Range("H138").Value = 2979.86
$ws.Range("I138").Value = 1314.8928

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 3156.138
$ws.Range("M32").Value = -2869.138
$ws.Range("H32").Value = 3166.1343
$ws.Range("I32").Value = 3156.138
$ws.Range("H45").Value = 3160.5
$ws.Range("I45").Value = 1451
$ws.Range("M45").Value = -1074
$ws.Range("K45").Value = 1451
$ws.Range("M60").Value = -65911.5
$ws.Range("K60").Value = 66644.5
$ws.Range("I60").Value = 66644.5
$ws.Range("H60").Value = 66644.5
$ws.Range("I74").Value = 894128.2
$ws.Range("H74").Value = 807984.5600000001
$ws.Range("M74").Value = -893254.2
$ws.Range("K74").Value = 894128.2
$ws.Range("M77").Value = -4466273
$ws.Range("K77").Value = 4470641
$ws.Range("I77").Value = 894128.2
$ws.Range("H77").Value = 807984.5600000001
$ws.Range("L98").Value = 130462
$ws.Range("J98").Value = 130462
$ws.Range("N98").Value = -136452
$ws.Range("H98").Value = 130462
$ws.Range("K110").Value = 2778.125
$ws.Range("I110").Value = 2778.125
$ws.Range("H110").Value = 3580.5557
$ws.Range("M110").Value = -733.125
$ws.Range("L135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L98").Value = 99999
$ws.Range("J98").Value = 99999
$ws.Range("N98").Value = -105989
$ws.Range("H98").Value = 99999
$ws.Range("H105").Value = 1906996.8
$ws.Range("J105").Value = 11466.333
$ws.Range("L105").Value = 11466.333
$ws.Range("N105").Value = -14960.333
$ws.Range("K107").Value = 1753.75
$ws.Range("I107").Value = 1753.75
$ws.Range("L107").Value = 3474.6667
$ws.Range("H107").Value = 2491.2856
$ws.Range("M107").Value = 166.25
$ws.Range("J107").Value = 3474.6667
$ws.Range("N107").Value = -7314.6667
$ws.Range("H132").Value = 200000
$ws.Range("J132").Value = 200000
$ws.Range("N132").Value = -210120
$ws.Range("L132").Value = 200000
$ws.Range("I134").Value = 2582.4546
$ws.Range("H134").Value = 8335700.5
$ws.Range("K134").Value = 7747.3638
$ws.Range("M134").Value = -5212.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36077816
$ws.Range("N31").Value = -780072.4
$ws.Range("L31").Value = 779482.4
$ws.Range("J31").Value = 779482.4
$ws.Range("L34").Value = 779482.4
$ws.Range("N34").Value = -779886.4
$ws.Range("H34").Value = 36077816
$ws.Range("J34").Value = 779482.4
$ws.Range("H42").Value = 5000
$ws.Range("K42").Value = 5000
$ws.Range("I42").Value = 5000
$ws.Range("M42").Value = -4407
$ws.Range("H102").Value = 86444.60000000001
$ws.Range("J102").Value = 86444.60000000001
$ws.Range("L102").Value = 86444.60000000001
$ws.Range("N102").Value = -91312.60000000001
$ws.Range("I134").Value = 3714.25
$ws.Range("H134").Value = 3878
$ws.Range("K134").Value = 11142.75
$ws.Range("M134").Value = -8607.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L34").Value = 47631.999
$ws.Range("N34").Value = -47799.999
$ws.Range("H34").Value = 12033
$ws.Range("J34").Value = 15877.333
$ws.Range("K55").Value = 7384.3842
$ws.Range("I55").Value = 2461.4614
$ws.Range("H55").Value = 5291.5
$ws.Range("M55").Value = -7207.3842
$ws.Range("K107").Value = 680.00001
$ws.Range("I107").Value = 226.66667
$ws.Range("L107").Value = 30338514
$ws.Range("H107").Value = 7584685
$ws.Range("M107").Value = 1239.99999
$ws.Range("J107").Value = 10112838
$ws.Range("N107").Value = -30342354
$ws.Range("H122").Value = 83750
$ws.Range("J122").Value = 1666.6666
$ws.Range("N122").Value = -19899.9994
$ws.Range("L122").Value = 14999.9994
$ws.Range("M138").Value = -18353.6
$ws.Range("K138").Value = 23493.6
$ws.Range("H138").Value = 8948.471
$ws.Range("I138").Value = 7831.2
$ws.Range("L140").Value = 33259.2
$ws.Range("H140").Value = 3720.2856
$ws.Range("J140").Value = 11086.4
$ws.Range("M140").Value = 924.875
$ws.Range("K140").Value = 4255.125
$ws.Range("N140").Value = -43619.2
$ws.Range("I140").Value = 1418.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("L18").Value = 0
$ws.Range("I113").Value = 4344.3335
$ws.Range("L113").Value = 18518518
$ws.Range("H113").Value = 4632888
$ws.Range("J113").Value = 18518518
$ws.Range("N113").Value = -18522858
$ws.Range("M113").Value = -2174.3335
$ws.Range("K113").Value = 4344.3335
$ws.Range("J139").Value = 152280.67
$ws.Range("N139").Value = -162560.67
$ws.Range("L139").Value = 152280.67
$ws.Range("H139").Value = 152280.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 24666.334
$ws.Range("J4").Value = 24666.334
$ws.Range("N4").Value = -24892.334
$ws.Range("L4").Value = 24666.334
$ws.Range("L28").Value = 24666.334
$ws.Range("J28").Value = 24666.334
$ws.Range("H28").Value = 24666.334
$ws.Range("N28").Value = -25130.334
$ws.Range("N37").Value = -24880.334
$ws.Range("H37").Value = 24666.334
$ws.Range("L37").Value = 24666.334
$ws.Range("J37").Value = 24666.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I29").Value = 10000
$ws.Range("K29").Value = 10000
$ws.Range("M29").Value = -9710
$ws.Range("H29").Value = 13799.6
$ws.Range("H33").Value = 19599
$ws.Range("L33").Value = 28331.666
$ws.Range("I33").Value = 6500
$ws.Range("N33").Value = -28831.666
$ws.Range("J33").Value = 28331.666
$ws.Range("M33").Value = -6250
$ws.Range("K33").Value = 6500
$ws.Range("L36").Value = 28331.666
$ws.Range("N36").Value = -28831.666
$ws.Range("H36").Value = 19599
$ws.Range("M36").Value = -6250
$ws.Range("J36").Value = 28331.666
$ws.Range("I36").Value = 6500
$ws.Range("K36").Value = 6500
$ws.Range("L101").Value = 17819
$ws.Range("J101").Value = 17819
$ws.Range("N101").Value = -24309
$ws.Range("H101").Value = 17819
$ws.Range("K107").Value = 4154.4375
$ws.Range("I107").Value = 1384.8125
$ws.Range("L107").Value = 11169.7062
$ws.Range("H107").Value = 2589.4546
$ws.Range("M107").Value = -2234.4375
$ws.Range("J107").Value = 3723.2354
$ws.Range("N107").Value = -15009.7062
$ws.Range("L113").Value = 2325
$ws.Range("H113").Value = 621.6667
$ws.Range("J113").Value = 775
$ws.Range("N113").Value = -6665
$ws.Range("I126").Value = 8301.08
$ws.Range("K126").Value = 24903.24
$ws.Range("M126").Value = -22433.24
$ws.Range("H126").Value = 7547.2
$ws.Range("J138").Value = 98884.5
$ws.Range("L138").Value = 98884.5
$ws.Range("N138").Value = -109164.5
$ws.Range("H138").Value = 83128.164
